$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "updated with strain names"
# Column F ("strain") was empty for most samples; fill in the strain name
# for each block of replicate rows. Rows 2-4 and 14-16 already held
# "KN99alpha" and stay untouched.
$strainByRows = @(
    @{ Rows = @(5, 6, 7);    Strain = "TDY1447" },
    @{ Rows = @(8, 9, 10);   Strain = "TYD1449" },
    @{ Rows = @(11, 12, 13); Strain = "TDY1448" },
    @{ Rows = @(17, 18, 19); Strain = "TDY1445" },
    @{ Rows = @(20, 21, 22); Strain = "TDY1441" },
    @{ Rows = @(23, 24, 25); Strain = "TDY1442" }
)

foreach ($group in $strainByRows) {
    foreach ($r in $group.Rows) {
        $ws.Range("F$r").Value = $group.Strain
    }
}

# Update the view to match where the author ended up: scrolled so row 14 is
# at the top, with F24:F25 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$win.Left = 15700

$ws.Range("F24:F25").Select()
